# Scheduled market-data refresh: update H:N (price/profit) columns
# across the ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets to the latest pull.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 153.5
$ws.Range("I9").Value = 145.5
$ws.Range("K9").Value = 145.5
$ws.Range("M9").Value = 23.5
# Row 15
$ws.Range("H15").Value = 539.0909
$ws.Range("I15").Value = 539.0909
$ws.Range("K15").Value = 1617.2727
$ws.Range("M15").Value = -1448.2727
# Row 17
$ws.Range("H17").Value = 2464.842
$ws.Range("J17").Value = 1996.1818
$ws.Range("L17").Value = 5988.5454
$ws.Range("N17").Value = -6324.5454
# Row 18
$ws.Range("H18").Value = 1445.9286
$ws.Range("I18").Value = 1445.9286
$ws.Range("K18").Value = 1445.9286
$ws.Range("M18").Value = -1161.9286
# Row 19
$ws.Range("H19").Value = 39.8
$ws.Range("J19").Value = 39
$ws.Range("L19").Value = 39
$ws.Range("N19").Value = -389
# Row 29
$ws.Range("H29").Value = 2644.5557
$ws.Range("I29").Value = 999
$ws.Range("J29").Value = 2850.25
$ws.Range("K29").Value = 2997
$ws.Range("L29").Value = 8550.75
$ws.Range("M29").Value = -2716
$ws.Range("N29").Value = -9112.75
# Row 31
$ws.Range("H31").Value = 92
$ws.Range("I31").Value = 92
$ws.Range("K31").Value = 276
$ws.Range("M31").Value = -46
# Row 40
$ws.Range("H40").Value = 2083.3333
$ws.Range("I40").Value = 1933.3334
$ws.Range("K40").Value = 1933.3334
$ws.Range("M40").Value = -1758.3334
# Row 64
$ws.Range("H64").Value = 4400
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 4400
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 4400
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -4896
# Row 67
$ws.Range("H67").Value = 4400
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 4400
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 4400
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -6116
# Row 70
$ws.Range("H70").Value = 112218.86
$ws.Range("J70").Value = 130713.664
$ws.Range("L70").Value = 392140.992
$ws.Range("N70").Value = -392680.992
# Row 73
$ws.Range("H73").Value = 112218.86
$ws.Range("J73").Value = 130713.664
$ws.Range("L73").Value = 392140.992
$ws.Range("N73").Value = -394012.992
# Row 74
$ws.Range("H74").Value = 6926
$ws.Range("I74").Value = 5000
$ws.Range("K74").Value = 5000
$ws.Range("M74").Value = -4064
# Row 77
$ws.Range("H77").Value = 6926
$ws.Range("I77").Value = 5000
$ws.Range("K77").Value = 25000
$ws.Range("M77").Value = -20320
# Row 111
$ws.Range("H111").Value = 2521.125
$ws.Range("I111").Value = 2652.7144
$ws.Range("K111").Value = 7958.1432
$ws.Range("M111").Value = -4891.1432
# Row 112
$ws.Range("H112").Value = 2345.4736
$ws.Range("J112").Value = 2610.375
$ws.Range("L112").Value = 7831.125
$ws.Range("N112").Value = -10047.125
# Row 116
$ws.Range("H116").Value = 7557.0835
$ws.Range("I116").Value = 6974.75
$ws.Range("K116").Value = 6974.75
$ws.Range("M116").Value = -3532.75
# Row 132
$ws.Range("H132").Value = 1066.3334
$ws.Range("I132").Value = 1111.9678
$ws.Range("J132").Value = 359
$ws.Range("K132").Value = 3335.9034
$ws.Range("L132").Value = 1077
$ws.Range("M132").Value = -805.9033999999997
$ws.Range("N132").Value = -6137
# Row 133
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
# Row 141
$ws.Range("H141").Value = 2371.3635
$ws.Range("I141").Value = 1408.5
$ws.Range("K141").Value = 4225.5
$ws.Range("M141").Value = 954.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 1997.875
$ws.Range("I45").Value = 1977.6
$ws.Range("K45").Value = 1977.6
$ws.Range("M45").Value = -1600.6
# Row 97
$ws.Range("H97").Value = 15000
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 15000
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 15000
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -15992

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 80
$ws.Range("H80").Value = 179.5
$ws.Range("I80").Value = 249.25
$ws.Range("K80").Value = 249.25
$ws.Range("M80").Value = 748.75
# Row 83
$ws.Range("H83").Value = 179.5
$ws.Range("I83").Value = 249.25
$ws.Range("K83").Value = 1246.25
$ws.Range("M83").Value = 3745.75
# Row 94
$ws.Range("H94").Value = 6100
$ws.Range("I94").Value = 6100
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 6100
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -5649
$ws.Range("N94").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 2180
$ws.Range("I16").Value = 950
$ws.Range("K16").Value = 950
$ws.Range("M16").Value = -663
# Row 31
$ws.Range("H31").Value = 4243.6113
$ws.Range("I31").Value = 1803.909
$ws.Range("J31").Value = 8077.4287
$ws.Range("K31").Value = 1803.909
$ws.Range("L31").Value = 8077.4287
$ws.Range("M31").Value = -1508.909
$ws.Range("N31").Value = -8667.4287
# Row 34
$ws.Range("H34").Value = 4243.6113
$ws.Range("I34").Value = 1803.909
$ws.Range("J34").Value = 8077.4287
$ws.Range("K34").Value = 1803.909
$ws.Range("L34").Value = 8077.4287
$ws.Range("M34").Value = -1601.909
$ws.Range("N34").Value = -8481.4287
# Row 58
$ws.Range("H58").Value = 3673.6667
$ws.Range("I58").Value = 1812.5555
$ws.Range("K58").Value = 1812.5555
$ws.Range("M58").Value = -1609.5555
# Row 99
$ws.Range("H99").Value = 13550.137
$ws.Range("I99").Value = 9684.1
$ws.Range("J99").Value = 16771.834
$ws.Range("K99").Value = 9684.1
$ws.Range("L99").Value = 16771.834
$ws.Range("M99").Value = -8186.1
$ws.Range("N99").Value = -19767.834
# Row 113
$ws.Range("H113").Value = 2180
$ws.Range("I113").Value = 950
$ws.Range("K113").Value = 950
$ws.Range("M113").Value = 1220
# Row 126
$ws.Range("H126").Value = 13550.137
$ws.Range("I126").Value = 9684.1
$ws.Range("J126").Value = 16771.834
$ws.Range("K126").Value = 29052.3
$ws.Range("L126").Value = 50315.50199999999
$ws.Range("M126").Value = -26582.3
$ws.Range("N126").Value = -55255.50199999999
# Row 134
$ws.Range("H134").Value = 3579.647
$ws.Range("I134").Value = 2829.2
$ws.Range("K134").Value = 8487.599999999999
$ws.Range("M134").Value = -5952.599999999999
# Row 136
$ws.Range("H136").Value = 3673.6667
$ws.Range("I136").Value = 1812.5555
$ws.Range("K136").Value = 5437.666499999999
$ws.Range("M136").Value = -2887.666499999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 66733.734
$ws.Range("I2").Value = 90940.82000000001
$ws.Range("J2").Value = 164.25
$ws.Range("K2").Value = 545644.92
$ws.Range("L2").Value = 985.5
$ws.Range("M2").Value = -545531.92
$ws.Range("N2").Value = -1211.5
# Row 33
$ws.Range("H33").Value = 5333.3335
$ws.Range("I33").Value = 5500
$ws.Range("J33").Value = 5000
$ws.Range("K33").Value = 33000
$ws.Range("L33").Value = 30000
$ws.Range("M33").Value = -32717
$ws.Range("N33").Value = -30566
# Row 55
$ws.Range("H55").Value = 79146.16
$ws.Range("I55").Value = 250475
$ws.Range("K55").Value = 751425
$ws.Range("M55").Value = -751248
# Row 97
$ws.Range("H97").Value = 60.4
$ws.Range("I97").Value = 67.333336
$ws.Range("J97").Value = 50
$ws.Range("K97").Value = 202.000008
$ws.Range("L97").Value = 150
$ws.Range("M97").Value = 293.999992
$ws.Range("N97").Value = -1142
# Row 139
$ws.Range("H139").Value = 3954.8125
$ws.Range("I139").Value = 2178.2
$ws.Range("J139").Value = 6915.8335
$ws.Range("K139").Value = 6534.599999999999
$ws.Range("L139").Value = 20747.5005
$ws.Range("M139").Value = -1394.599999999999
$ws.Range("N139").Value = -31027.5005

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2230.7693
$ws.Range("I7").Value = 2091.182
$ws.Range("K7").Value = 2091.182
$ws.Range("M7").Value = -1979.182
# Row 126
$ws.Range("H126").Value = 2230.7693
$ws.Range("I126").Value = 2091.182
$ws.Range("K126").Value = 6273.545999999999
$ws.Range("M126").Value = -3803.545999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 952.8333
$ws.Range("I107").Value = 422.83334
$ws.Range("K107").Value = 1268.50002
$ws.Range("M107").Value = 651.4999800000001
# Row 113
$ws.Range("H113").Value = 844.8
$ws.Range("I113").Value = 743.4545000000001
$ws.Range("J113").Value = 1123.5
$ws.Range("K113").Value = 2230.3635
$ws.Range("L113").Value = 3370.5
$ws.Range("M113").Value = -60.36350000000039
$ws.Range("N113").Value = -7710.5
# Row 126
$ws.Range("H126").Value = 2904.2856
$ws.Range("I126").Value = 816
$ws.Range("K126").Value = 2448
$ws.Range("M126").Value = 22
# Row 132
$ws.Range("H132").Value = 2694.4
$ws.Range("I132").Value = 2118
$ws.Range("K132").Value = 6354
$ws.Range("M132").Value = -3824
# Row 136
$ws.Range("H136").Value = 2804.818
$ws.Range("I136").Value = 1068.5
$ws.Range("K136").Value = 3205.5
$ws.Range("M136").Value = -655.5
